$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: convert politeness_score to a real number and clear the
# polite_expressions ("nan") placeholder text.
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = ""

# Row 13: new annotation row appended for parisk.
$ws.Range("A13").Value = "parisk"
$ws.Range("B13").Value = "2"
$ws.Range("C13").Value = "nan"
$ws.Range("D13").Value = "DIS"
$ws.Range("E13").Value = "WRI"
$ws.Range("F13").Value = "77474e59-42ef-43e4-850b-a07d6b41a266"
$ws.Range("G13").Value = "Syg-YfWCW_annotated.xlsx"
$ws.Range("H13").Value = "You absolutely know this but you hide these results."
